$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New journal entries ("Remplissage du journal de bord")
$textB23 = "J'ai listé les articles dans la page product.php. Tous les articles que j'ai dans ma base de données sont affichés dans cette page. J'ai aussi mis un filtre sur les articles pour avoir le choix d'afficher uniquement des habits, des chaussures ou seulement les sacs à dos, mais j'eu eu un souci avec cette fonction."
$textB24 = "J'ai rélgé mon problème de fonction que j'avais la veille. Il est maintenant possible d'appliquer les filtres sur les types d'articles que l'on souhaite. "

# Row 23 - copy formatting from row 21 (date cell style, wrapped-text style, plain style)
# so that the existing number-format / wrap-text styles are reused instead of new ones being minted.
$ws.Range("A21").Copy($ws.Range("A23"))
$ws.Range("A23").Value = 43146
$ws.Range("B21").Copy($ws.Range("B23"))
$ws.Range("B23").Value = $textB23
$ws.Range("C21").Copy($ws.Range("C23"))
$ws.Range("C23").Value = "4 périodes"
$ws.Rows.Item(23).RowHeight = 60

# Row 24
$ws.Range("A22").Copy($ws.Range("A24"))
$ws.Range("A24").Value = 43147
$ws.Range("B22").Copy($ws.Range("B24"))
$ws.Range("B24").Value = $textB24
$ws.Range("C17").Copy($ws.Range("C24"))
$ws.Range("C24").Value = "1 période"
$ws.Rows.Item(24).RowHeight = 30

# Update selection to reflect the new extent of the journal
$ws.Range("C25").Select()
